$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting,
# e.g. trailing zeros / exact decimal digits), by forcing Text number format
# before assigning the value, for the cells whose new value looks like a plain number.

$ws.Range("D2").Value = '23.646.35'
$ws.Range("E2").Value = '  +1.83%  '
$ws.Range("D3").Value = '1.655.12'
$ws.Range("E3").Value = '  +3.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9954'
$ws.Range("E4").Value = '  -0.57%  '
$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9965'
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.38'
$ws.Range("E6").Value = '  +0.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3765'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.95'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3650'
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.266'
$ws.Range("E10").Value = '  -0.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08148'
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9956'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.17'
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.680'
$ws.Range("E14").Value = '  +1.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001277'
$ws.Range("E15").Value = '  +2.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.370'
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").Value = '1.648.32'
$ws.Range("E17").Value = '  +3.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.05'
$ws.Range("E18").Value = '  +1.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06882'
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.35'
$ws.Range("E20").Value = '  +1.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.593'
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9974'
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("D23").Value = '23.641.45'
$ws.Range("E23").Value = '  +1.84%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.171'
$ws.Range("E25").Value = '  +4.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.408'
$ws.Range("E26").Value = '  -1.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.39'
$ws.Range("E27").Value = '  +0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.91'
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.308'
$ws.Range("E29").Value = '  +0.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.71'
$ws.Range("E30").Value = '  +0.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.302'
$ws.Range("E31").Value = '  -3.20%  '
$ws.Range("D32").Value = '1.831.65'
$ws.Range("E32").Value = '  +3.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.873'
$ws.Range("E33").Value = '  +2.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9722'
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '10.74'
$ws.Range("E35").Value = '  +3.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02830'
$ws.Range("E36").Value = '  +3.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.320'
$ws.Range("E37").Value = '  +3.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.07376'
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2551'
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08879'
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.379'
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7137'
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.44'
$ws.Range("E43").Value = '  +5.08%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.62'
$ws.Range("E44").Value = '  +1.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6591'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.359'
$ws.Range("E46").Value = '  +1.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9957'
$ws.Range("E47").Value = '  -0.47%  '
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08036'
$ws.Range("E49").Value = '  +1.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '129.58'
$ws.Range("E50").Value = '  -2.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.219'
$ws.Range("E51").Value = '  +1.45%  '
